$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 23750
$ws.Range("D2").Value = -0.0084
$ws.Range("I2").Value = 4.21
$ws.Range("J2").Value = 82
$ws.Range("K2").Value = 82

$ws.Range("C3").Value = 101800
$ws.Range("D3").NumberFormat = "0%"
$ws.Range("D3").Value = 0
$ws.Range("I3").Value = 6.39
$ws.Range("J3").Value = 65
$ws.Range("K3").Value = 65

$ws.Range("C4").Value = 443000
$ws.Range("D4").Value = 0.0126
$ws.Range("I4").Value = 4.29
$ws.Range("J4").Value = 77
$ws.Range("K4").Value = 77

$ws.Range("C5").Value = 30350
$ws.Range("D5").Value = 0.0117
$ws.Range("I5").Value = 6.59
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 45

$ws.Range("C6").Value = 29650
$ws.Range("D6").Value = -0.0279
$ws.Range("I6").Value = 4.05
$ws.Range("J6").Value = 75
$ws.Range("K6").Value = 75

$ws.Range("C7").Value = 25100
$ws.Range("D7").Value = 0.0101
$ws.Range("I7").Value = 4.78
$ws.Range("J7").Value = 68
$ws.Range("K7").Value = 68

$ws.Range("C8").Value = 10760
$ws.Range("D8").Value = 0.0132
$ws.Range("I8").Value = 4.79
$ws.Range("J8").Value = 87
$ws.Range("K8").Value = 87

$ws.Range("C9").Value = 87500
$ws.Range("D9").Value = 0.0139
$ws.Range("I9").Value = 3.43
$ws.Range("J9").Value = 81
$ws.Range("K9").Value = 81

$ws.Range("C10").Value = 212500
$ws.Range("D10").Value = 0.0095
$ws.Range("I10").Value = 5.65
$ws.Range("J10").Value = 42
$ws.Range("K10").Value = 42

$ws.Range("C11").Value = 129400
$ws.Range("D11").Value = 0.0094
$ws.Range("I11").Value = 5.26
$ws.Range("J11").Value = 83
$ws.Range("K11").Value = 83

$ws.Range("C12").Value = 19600
$ws.Range("D12").Value = 0.0015
$ws.Range("I12").Value = 4.85
$ws.Range("J12").Value = 74
$ws.Range("K12").Value = 74

$ws.Range("C13").Value = 70600
$ws.Range("D13").Value = -0.0028
$ws.Range("I13").Value = 4.96
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 80

$ws.Range("C15").Value = 82700
$ws.Range("D15").Value = 0.0122
$ws.Range("I15").Value = 6.65
$ws.Range("J15").Value = 87
$ws.Range("K15").Value = 87

$ws.Range("C16").Value = 19370
$ws.Range("D16").Value = 0.0016
$ws.Range("I16").Value = 5.5
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 79

$ws.Range("C17").Value = 50400
$ws.Range("D17").NumberFormat = "0%"
$ws.Range("D17").Value = 0
$ws.Range("I17").Value = 5.56
$ws.Range("J17").Value = 71
$ws.Range("K17").Value = 71

$ws.Range("C18").Value = 20050
$ws.Range("D18").Value = 0.0106
$ws.Range("I18").Value = 6.13
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = 36

$ws.Range("C19").Value = 54300
$ws.Range("D19").Value = -0.0018
$ws.Range("I19").Value = 3.68
$ws.Range("J19").Value = 87
$ws.Range("K19").Value = 87

$ws.Range("C20").Value = 14560
$ws.Range("D20").Value = -0.0082
$ws.Range("I20").Value = 4.46
$ws.Range("J20").Value = 75
$ws.Range("K20").Value = 75

$ws.Range("C21").Value = 138500
$ws.Range("D21").Value = -0.0142
$ws.Range("I21").Value = 3.9
$ws.Range("J21").Value = 92
$ws.Range("K21").Value = 92

$ws.Range("C22").Value = 42300
$ws.Range("D22").Value = -0.0047
$ws.Range("I22").Value = 3.44
$ws.Range("J22").Value = 47
$ws.Range("K22").Value = 47

$ws.Range("C23").Value = 69200
$ws.Range("D23").Value = 0.0073
$ws.Range("I23").Value = 3.12
$ws.Range("J23").Value = 91
$ws.Range("K23").Value = 91

$ws.Range("C24").Value = 48200
$ws.Range("D24").Value = -0.0031
$ws.Range("I24").Value = 5.6
$ws.Range("J24").Value = 66
$ws.Range("K24").Value = 66

$ws.Range("C25").Value = 85800
$ws.Range("D25").Value = 0.0035
$ws.Range("I25").Value = 4.2
$ws.Range("J25").Value = 84
$ws.Range("K25").Value = 84

$ws.Range("C26").Value = 115000
$ws.Range("D26").Value = -0.0052
$ws.Range("I26").Value = 2.76
$ws.Range("J26").Value = 87
$ws.Range("K26").Value = 87

$ws.Range("C27").Value = 14620
$ws.Range("D27").Value = 0.0069
$ws.Range("I27").Value = 4.45
$ws.Range("J27").Value = 87
$ws.Range("K27").Value = 87

$ws.Range("C28").Value = 14050
$ws.Range("D28").Value = 0.0086
$ws.Range("I28").Value = 3.56
$ws.Range("J28").Value = 85
$ws.Range("K28").Value = 85

$ws.Range("C29").Value = 22500
$ws.Range("D29").Value = -0.0044
$ws.Range("I29").Value = 4.42
$ws.Range("J29").Value = 83
$ws.Range("K29").Value = 83

$ws.Range("C30").Value = 25350
$ws.Range("D30").Value = 0.012
$ws.Range("I30").Value = 4.73
$ws.Range("J30").Value = 91
$ws.Range("K30").Value = 91

$ws.Range("D14").Value = -0.0142

$ws.Range("H17").Select()
